# Edit ILLINOIS_2019.xlsx per the commit diff:
#  1. Rename header row (A1:D1) to snake_case column names.
#  2. Title-case the Spanish connector words (de/del/la/las/los/el/y) inside
#     every state/municipality name in columns A and B (rows 2..1866).
#  3. Apply the tiny (1-ULP) floating point corrections to the few D cells
#     that the source recalculation produced.
#  4. Remove the trailing footnote rows (1868..1872) and the blank row 1867,
#     shrinking the sheet's dimension down to A1:D1866.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header renames ------------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Title-case connector words in columns A and B, rows 2..1866 --------
# NOTE: this runtime's `-eq`/`-ne` string comparisons are case-INSENSITIVE
# (even with the `-c` prefix), so we don't rely on equality checks here -
# the cell is simply rewritten with the transformed text every time (a
# no-op write when there was nothing to change).
for ($r = 2; $r -le 1866; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $valA = $cellA.Value()
    if ($valA -ne $null -and $valA -ne "") {
        $newA = $valA -replace '\bde\b', 'De'
        $newA = $newA -replace '\bdel\b', 'Del'
        $newA = $newA -replace '\blas\b', 'Las'
        $newA = $newA -replace '\blos\b', 'Los'
        $newA = $newA -replace '\bla\b', 'La'
        $newA = $newA -replace '\bel\b', 'El'
        $newA = $newA -replace '\by\b', 'Y'
        $cellA.Value = $newA
    }

    $cellB = $ws.Cells.Item($r, 2)
    $valB = $cellB.Value()
    if ($valB -ne $null -and $valB -ne "") {
        $newB = $valB -replace '\bde\b', 'De'
        $newB = $newB -replace '\bdel\b', 'Del'
        $newB = $newB -replace '\blas\b', 'Las'
        $newB = $newB -replace '\blos\b', 'Los'
        $newB = $newB -replace '\bla\b', 'La'
        $newB = $newB -replace '\bel\b', 'El'
        $newB = $newB -replace '\by\b', 'Y'
        $cellB.Value = $newB
    }
}

# --- 3. Floating point (1-ULP) corrections on D column ---------------------
$ws.Range("D414").Value = 0.0009114055899542852
$ws.Range("D423").Value = 0.0009548058561425844
$ws.Range("D625").Value = 0.0009114055899542852
$ws.Range("D846").Value = 0.0009114055899542852
$ws.Range("D1201").Value = 0.0009548058561425844
$ws.Range("D1785").Value = 0.0009114055899542852

# --- 4. Drop the trailing footnote rows and shrink the sheet ---------------
$ws.Rows("1868:1872").Delete()
